$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.34837547187382256
$ws.Range("B1").Value = 0.34742996497635659
$ws.Range("A2").Value = -0.27242477466130488
$ws.Range("B2").Value = 0.26920617902004551
$ws.Range("A3").Value = -0.16625572935188515
$ws.Range("B3").Value = 0.16526299472259609
$ws.Range("A4").Value = -0.15326299496470952
$ws.Range("B4").Value = 0.1523697339301826
$ws.Range("A5").Value = -0.14636973488085481
$ws.Range("B5").Value = 0.14457400213812743
$ws.Range("A6").Value = -0.027736976715153183
$ws.Range("B6").Value = 0.027726301844821943
$ws.Range("A7").Value = -0.0077263029902017166
$ws.Range("B7").Value = 0.0077253221862001453
$ws.Range("A8").Value = 0.012274676667834328
$ws.Range("B8").Value = -0.012291800861471636
$ws.Range("A9").Value = 0.01829179988336449
$ws.Range("B9").Value = -0.018325222740815938
$ws.Range("A10").Value = 0.024325221767050209
$ws.Range("B10").Value = -0.024325469848236025
$ws.Range("A11").Value = -0.022212180006071947
$ws.Range("B11").Value = 0.022196694234096981
$ws.Range("A12").Value = -0.016196695208564371
$ws.Range("B12").Value = 0.016150714835222679
$ws.Range("A13").Value = -0.010150715814274847
$ws.Range("B13").Value = 0.010143316343850017
$ws.Range("A14").Value = -0.052570384720074692
$ws.Range("B14").Value = 0.052386003378523505
$ws.Range("A15").Value = -0.046386004365240652
$ws.Range("B15").Value = 0.046187039899884041
$ws.Range("A16").Value = -0.015027240807452991
$ws.Range("B16").Value = 0.01500389052542328
$ws.Range("A17").Value = -0.0090038915238270789
$ws.Range("B17").Value = 0.0089999989653213319
$ws.Range("A18").Value = -0.036106070532795798
$ws.Range("B18").Value = 0.036096115311814003
$ws.Range("A19").Value = -0.027096116244976542
$ws.Range("B19").Value = 0.027012896569189593
$ws.Range("A20").Value = -0.018012897510670811
$ws.Range("B20").Value = 0.018004208142052391
$ws.Range("A21").Value = -0.0090042090847557432
$ws.Range("B21").Value = 0.0089999990564555432
$ws.Range("A22").Value = -0.093950450629849414
$ws.Range("B22").Value = 0.093637175674693651
$ws.Range("A23").Value = -0.084637176645047774
$ws.Range("B23").Value = 0.084127099234314606
$ws.Range("A24").Value = -0.042127100607114443
$ws.Range("B24").Value = 0.041999998619611567
$ws.Range("A25").Value = -0.094963624307755623
$ws.Range("B25").Value = 0.094719794909021715
$ws.Range("A26").Value = -0.088719795892107101
$ws.Range("B26").Value = 0.088402714021366791
$ws.Range("A27").Value = -0.082402715010000627
$ws.Range("B27").Value = 0.08130685329541576
$ws.Range("A28").Value = -0.075306854305498661
$ws.Range("B28").Value = 0.074547835032388932
$ws.Range("A29").Value = -0.062547836127668788
$ws.Range("B29").Value = 0.062172878258337505
$ws.Range("A30").Value = -0.042172879455150802
$ws.Range("B30").Value = 0.042018753206073622
$ws.Range("A31").Value = -0.02701875435392509
$ws.Range("B31").Value = 0.027000401454360556
$ws.Range("A32").Value = -0.0060004026742142003
$ws.Range("B32").Value = 0.0059999989566916767
